$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 2008 and 2009 rows (rows 2 and 3) are removed from the table;
# all subsequent rows (2010..2019) shift up by two rows, and the
# sheet's used range shrinks from A1:E13 to A1:E11.
$ws.Rows("2:3").Delete()
